# correção nos dados e inicio da analise PNAD 2009
# The "unnamed: 1_level_1" and "unnamed: 5_level_1" placeholder labels
# (auto-generated by pandas for sub-header cells with no real name) are
# corrected to "total", matching the sibling sub-header already in C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
